$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Quantities executed upto date (column C) - plain numeric cells ---
$ws.Range("C8").Value = 7
$ws.Range("C9").Value = 62
$ws.Range("C10").Value = 11
$ws.Range("C11").Value = 25
$ws.Range("C12").Value = 76
$ws.Range("C13").Value = 17
$ws.Range("C14").Value = 10
$ws.Range("C15").Value = 71
$ws.Range("C16").Value = 55
$ws.Range("C17").Value = 37

# --- Upto date Amount (column G / H) - these are text cells formatted
# like "12345.00" (number stored as text), so we force text entry and
# then restore the original "Normal" cell style so no formatting drifts. ---
$textCells = @("G9", "G10", "G11", "G13", "G14", "G19", "H19", "G21", "H21")
$textValues = @("15872.00", "5192.00", "16550.00", "2312.00", "230.00", "40156.00", "40156.00", "40156.00", "40156.00")

for ($i = 0; $i -lt $textCells.Count; $i++) {
    $cell = $ws.Range($textCells[$i])
    $cell.NumberFormat = "@"
    $cell.Value = $textValues[$i]
    $cell.Style = "Normal"
}
